$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (column D) updates ---
$ws.Range("D2").Value = '60.620.05'
$ws.Range("D3").Value = '2.409.12'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '564.22'
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '137.57'
$ws.Range("D6").Style = "Normal"
$ws.Range("D9").Value = '2.392.23'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.62'
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = '2.823.08'
$ws.Range("D17").Value = '60.655.97'
$ws.Range("D18").Value = '2.407.50'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '8.01'
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '322.53'
$ws.Range("D21").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.81'
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '64.10'
$ws.Range("D26").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '553.82'
$ws.Range("D28").Style = "Normal"
$ws.Range("D30").Value = '0.0₃0919'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.30'
$ws.Range("D32").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '152.37'
$ws.Range("D37").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.53'
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.11'
$ws.Range("D40").Style = "Normal"
$ws.Range("D45").Value = '0.0₆0293'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '141.94'
$ws.Range("D46").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.582'
$ws.Range("D48").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '19.08'
$ws.Range("D50").Style = "Normal"

# --- Volume(1h) (column E) updates ---
$ws.Range("E2").Value = '  -1.80%  '
$ws.Range("E3").Value = '  -1.44%  '
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("E5").Value = '  -2.62%  '
$ws.Range("E6").Value = '  -2.60%  '
$ws.Range("E7").Value = '  +0.19%  '
$ws.Range("E8").Value = '  +1.29%  '
$ws.Range("E9").Value = '  -1.78%  '
$ws.Range("E10").Value = '  -3.31%  '
$ws.Range("E11").Value = '  -0.71%  '
$ws.Range("E12").Value = '  -2.53%  '
$ws.Range("E13").Value = '  -1.11%  '
$ws.Range("E14").Value = '  -0.67%  '
$ws.Range("E15").Value = '  -2.31%  '
$ws.Range("E16").Value = '  -2.32%  '
$ws.Range("E17").Value = '  -1.61%  '
$ws.Range("E18").Value = '  -1.24%  '
$ws.Range("E19").Value = '  +11.98%  '
$ws.Range("E20").Value = '  -0.84%  '
$ws.Range("E21").Value = '  -0.50%  '
$ws.Range("E22").Value = '  -0.85%  '
$ws.Range("E23").Value = '  -5.86%  '
$ws.Range("E24").Value = '  +0.03%  '
$ws.Range("E25").Value = '  -6.53%  '
$ws.Range("E26").Value = '  -1.21%  '
$ws.Range("E27").Value = '  -9.74%  '
$ws.Range("E28").Value = '  -4.45%  '
$ws.Range("E30").Value = '  -0.55%  '
$ws.Range("E32").Value = '  -5.25%  '
$ws.Range("E33").Value = '  -3.91%  '
$ws.Range("E34").Value = '  -1.17%  '
$ws.Range("E35").Value = '  +0.21%  '
$ws.Range("E36").Value = '  +1.12%  '
$ws.Range("E37").Value = '  +0.38%  '
$ws.Range("E38").Value = '  -1.26%  '
$ws.Range("E39").Value = '  -4.79%  '
$ws.Range("E40").Value = '  -0.63%  '
$ws.Range("E41").Value = '  -1.18%  '
$ws.Range("E42").Value = '  +0.04%  '
$ws.Range("E43").Value = '  -1.43%  '
$ws.Range("E44").Value = '  -0.20%  '
$ws.Range("E45").Value = '  +5.86%  '
$ws.Range("E46").Value = '  +0.83%  '
$ws.Range("E47").Value = '  -1.72%  '
$ws.Range("E48").Value = '  -2.42%  '
$ws.Range("E49").Value = '  -1.89%  '
$ws.Range("E50").Value = '  -2.56%  '
$ws.Range("E51").Value = '  -0.08%  '
